$fundData = @(
    @('513330','华夏恒生互联网科技业ETF（QDII）','233.65','96.98','4.51','10.5376',8),
    @('513180','华夏恒生科技交易型开放式指数证券投资基金（QDII）','89.43','94.63','4.87','4.3552',9),
    @('011401','汇添富成长精选混合A','42.03','87.99','3.73','1.5677',8),
    @('011136','广发盛兴混合A','22.19','93.21','6.50','1.4424',6),
    @('513010','易方达恒生科技交易型开放式指数证券投资基金（QDII）','26.24','94.46','4.83','1.2674',9),
    @('010792','华安成长先锋混合A','16.44','93.43','7.07','1.1623',2),
    @('010161','广发瑞安精选股票A','8.03','92.30','6.44','0.5171',5),
    @('501311','嘉实恒生港股通新经济指数（LOF）A','22.36','94.24','2.28','0.5098',10),
    @('009362','招商丰盈积极配置混合A','21.39','87.04','1.72','0.3679',10),
    @('159740','大成恒生科技交易型开放式指数证券投资基金（QDII）','7.06','98.90','5.05','0.3565',9),
    @('014887','招商安福1年定期开放债券','17.22','27.65','2.05','0.3530',5),
    @('012208','华夏港股前沿经济混合型证券投资基金（QDII）A','12.53','82.34','2.73','0.3421',5),
    @('010793','华安成长先锋混合C','3.80','93.43','7.07','0.2687',2),
    @('513580','华安恒生科技交易型开放式指数证券投资基金（QDII）','4.91','93.60','4.89','0.2401',9),
    @('513980','景顺长城中证港股通科技交易型开放式指数证券投资基金','5.03','97.36','4.02','0.2022',8),
    @('006614','嘉实恒生港股通新经济指数（LOF）C','8.52','94.24','2.28','0.1943',10),
    @('159742','博时恒生科技交易型开放式指数证券投资基金(QDII)','3.73','96.57','4.94','0.1843',9),
    @('159741','嘉实恒生科技交易型开放式指数证券投资基金（QDII）','3.29','99.73','5.10','0.1678',9),
    @('006122','华安低碳生活混合','4.43','86.68','3.70','0.1639',3),
    @('513860','海富通中证港股通科技交易型开放式指数证券投资基金','3.87','95.13','4.15','0.1606',8),
    @('012010','富国泰享回报6个月持有期混合型证券投资基金A','9.29','29.91','1.50','0.1394',2),
    @('011137','广发盛兴混合C','2.10','93.21','6.50','0.1365',6),
    @('513890','上投摩根恒生科技ETF（QDII）','2.10','93.73','4.80','0.1008',9),
    @('159747','南方中证香港科技交易型开放式指数证券投资基金(QDII)','3.07','100.03','2.95','0.0906',10),
    @('009695','招商成长精选一年定期开放混合A','5.11','87.04','1.49','0.0761',10),
    @('013127','汇添富恒生科技指数（QDII）A','1.51','91.71','4.68','0.0707',9),
    @('009363','招商丰盈积极配置混合C','4.09','87.04','1.72','0.0703',10),
    @('011402','汇添富成长精选混合C','1.65','87.99','3.73','0.0615',8),
    @('513160','银华恒生港股通中国科技ETF','0.62','92.07','9.03','0.0560',4),
    @('003993','前海开源沪港深核心驱动灵活配置混合','0.58','82.10','8.43','0.0489',3),
    @('010162','广发瑞安精选股票C','0.65','92.30','6.44','0.0419',5),
    @('159751','鹏华中证港股通科技ETF','0.90','91.30','3.86','0.0347',8),
    @('011144','华安汇宏精选混合A','1.07','85.87','3.05','0.0326',7),
    @('013128','汇添富恒生科技指数（QDII）C','0.63','91.71','4.68','0.0295',9),
    @('513020','国泰中证港股通科技ETF','0.76','91.35','3.52','0.0268',8),
    @('517360','华安中证沪港深科技100交易型开放式指数证券投资基金','0.72','93.71','3.34','0.0240',8),
    @('009696','招商成长精选一年定期开放混合C','1.36','87.04','1.49','0.0203',10),
    @('012209','华夏港股前沿经济混合型证券投资基金（QDII）C','0.46','82.34','2.73','0.0126',5),
    @('011145','华安汇宏精选混合C','0.25','85.87','3.05','0.0076',7),
    @('001900','诺安精选价值混合','0.13','89.96','3.62','0.0047',6),
    @('006477','中邮沪港深精选混合','0.05','83.91','6.78','0.0034',6),
    @('012011','富国泰享回报6个月持有期混合型证券投资基金C','0.09','29.91','1.50','0.0014',2),
    @('004321','前海开源沪港深强国产业灵活配置混合','0.02','64.32','4.65','0.0009',2),
)

$wb = $excel.ActiveWorkbook

# --- Step 1: Insert the new "2022-Q1" sheet after "2021-Q4", before "总计" ---
$prevSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $prevSheet)
$newSheet.Name = "2022-Q1"

# Use the existing "2021-Q4" sheet as a format template (same column layout).
$template = $wb.Worksheets.Item("2021-Q4")

# Copy header-row formatting (bold/border/center) for B1:H1
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Copy column-A formatting (bold/border/center) down through row 44
$template.Range("A2").Copy()
$newSheet.Range("A2:A44").PasteSpecial(-4122)

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (row 2..44)
for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $row = $fundData[$i]

    $newSheet.Cells.Item($r, 1).Value = $i

    # Columns B..G are stored as TEXT even though several look numeric.
    for ($col = 2; $col -le 7; $col++) {
        $cell = $newSheet.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col - 2]
        $cell.Style = "Normal"
    }

    # Column H (仓位排名) is a real number.
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# --- Step 2: Prepend a "2022-Q1" row to the "总计" summary sheet ---
$total = $wb.Worksheets.Item("总计")

# Snapshot the 5 existing data rows (rows 2-6) before shifting them down.
$existing = @()
for ($r = 2; $r -le 6; $r++) {
    $existing += ,@($total.Cells.Item($r, 2).Value2, $total.Cells.Item($r, 3).Value2, $total.Cells.Item($r, 4).Value2)
}

# Copy row-2 formatting down onto the new row 7 (which had no prior formatting).
$total.Range("A2:D2").Copy()
$total.Range("A7:D7").PasteSpecial(-4122)

# Re-write rows 3..7 with the snapshotted values (shifted down by one row).
for ($i = 0; $i -lt $existing.Count; $i++) {
    $r = $i + 3
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $existing[$i][0]
    $total.Cells.Item($r, 3).Value = $existing[$i][1]
    $total.Cells.Item($r, 4).Value = $existing[$i][2]
}

# New first data row: 2022-Q1 summary.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 43
$total.Cells.Item(2, 4).Value = 25.45
